$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 63.91118233333333
$ws.Range("H2").Value = 191.733547
$ws.Range("I2").Value = 0.4067926910433548
$ws.Range("J2").Value = 0.4067926910433549
$ws.Range("M2").Value = 1.164012
$ws.Range("N2").Value = 3.492036
$ws.Range("O2").Value = 0.02222380689314669
$ws.Range("P2").Value = 0.02222380689314669
$ws.Range("Q2").Value = 74.39338317018799
$ws.Range("R2").Value = 669.5404485316919
$ws.Range("S2").Value = 0.009040482211291
$ws.Range("T2").Value = 0.009040482211291001

$ws.Range("G3").Value = 63.91118233333333
$ws.Range("H3").Value = 191.733547
$ws.Range("I3").Value = 0.4067926910433548
$ws.Range("J3").Value = 0.4067926910433549
$ws.Range("M3").Value = 14.70158366666666
$ws.Range("N3").Value = 44.10475099999999
$ws.Range("O3").Value = 0.2806888214480945
$ws.Range("P3").Value = 0.2806888214480945
$ws.Range("Q3").Value = 939.5955943090884
$ws.Range("R3").Value = 8456.360348781795
$ws.Range("S3").Value = 0.1141821610226581
$ws.Range("T3").Value = 0.1141821610226581

$ws.Range("G4").Value = 63.91118233333333
$ws.Range("H4").Value = 191.733547
$ws.Range("I4").Value = 0.4067926910433548
$ws.Range("J4").Value = 0.4067926910433549
$ws.Range("M4").Value = 36.51120933333333
$ws.Range("N4").Value = 109.533628
$ws.Range("O4").Value = 0.6970873716587588
$ws.Range("P4").Value = 0.6970873716587588
$ws.Range("Q4").Value = 2333.474556913168
$ws.Range("R4").Value = 21001.27101221852
$ws.Range("S4").Value = 0.2835700478094057
$ws.Range("T4").Value = 0.2835700478094058

$ws.Range("I5").Value = 0.3656254573230189
$ws.Range("J5").Value = 0.365625457323019
$ws.Range("M5").Value = 1.164012
$ws.Range("N5").Value = 3.492036
$ws.Range("O5").Value = 0.02222380689314669
$ws.Range("P5").Value = 0.02222380689314669
$ws.Range("Q5").Value = 66.86480692079999
$ws.Range("R5").Value = 601.7832622871999
$ws.Range("S5").Value = 0.008125589558765219
$ws.Range("T5").Value = 0.008125589558765219

$ws.Range("I6").Value = 0.3656254573230189
$ws.Range("J6").Value = 0.365625457323019
$ws.Range("M6").Value = 14.70158366666666
$ws.Range("N6").Value = 44.10475099999999
$ws.Range("O6").Value = 0.2806888214480945
$ws.Range("P6").Value = 0.2806888214480945
$ws.Range("Q6").Value = 844.5089511977998
$ws.Range("R6").Value = 7600.580560780199
$ws.Range("S6").Value = 0.1026269787074188
$ws.Range("T6").Value = 0.1026269787074188

$ws.Range("I7").Value = 0.3656254573230189
$ws.Range("J7").Value = 0.365625457323019
$ws.Range("M7").Value = 36.51120933333333
$ws.Range("N7").Value = 109.533628
$ws.Range("O7").Value = 0.6970873716587588
$ws.Range("P7").Value = 0.6970873716587588
$ws.Range("Q7").Value = 2097.3280022184
$ws.Range("R7").Value = 18875.9520199656
$ws.Range("S7").Value = 0.254872889056835
$ws.Range("T7").Value = 0.254872889056835

$ws.Range("G8").Value = 35.755375
$ws.Range("H8").Value = 107.266125
$ws.Range("I8").Value = 0.2275818516336261
$ws.Range("J8").Value = 0.2275818516336262
$ws.Range("M8").Value = 1.164012
$ws.Range("N8").Value = 3.492036
$ws.Range("O8").Value = 0.02222380689314669
$ws.Range("P8").Value = 0.02222380689314669
$ws.Range("Q8").Value = 41.6196855645
$ws.Range("R8").Value = 374.5771700805
$ws.Range("S8").Value = 0.005057735123090467
$ws.Range("T8").Value = 0.005057735123090467

$ws.Range("G9").Value = 35.755375
$ws.Range("H9").Value = 107.266125
$ws.Range("I9").Value = 0.2275818516336261
$ws.Range("J9").Value = 0.2275818516336262
$ws.Range("M9").Value = 14.70158366666666
$ws.Range("N9").Value = 44.10475099999999
$ws.Range("O9").Value = 0.2806888214480945
$ws.Range("P9").Value = 0.2806888214480945
$ws.Range("Q9").Value = 525.6606370955416
$ws.Range("R9").Value = 4730.945733859874
$ws.Range("S9").Value = 0.06387968171801763
$ws.Range("T9").Value = 0.06387968171801764

$ws.Range("G10").Value = 35.755375
$ws.Range("H10").Value = 107.266125
$ws.Range("I10").Value = 0.2275818516336261
$ws.Range("J10").Value = 0.2275818516336262
$ws.Range("M10").Value = 36.51120933333333
$ws.Range("N10").Value = 109.533628
$ws.Range("O10").Value = 0.6970873716587588
$ws.Range("P10").Value = 0.6970873716587588
$ws.Range("Q10").Value = 1305.471981416833
$ws.Range("R10").Value = 11749.2478327515
$ws.Range("S10").Value = 0.158644434792518
$ws.Range("T10").Value = 0.1586444347925181
